{"js": "// Replace the placeholder path text \"{path to org-opensim-rcnl.jar}\" with\n// \"[PATH_TO_NMSM-CORE]\\gui\\org-opensim-rcnl.jar\" in the GUI installation\n// instructions paragraph (\"Run the command \"opensim64 --reload ...\").\nconst results = context.document.body.search(\"{path to org-opensim-rcnl.jar}\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text not found\");\n}\n\nresults.items[0].insertText(\"[PATH_TO_NMSM-CORE]\\\\gui\\\\org-opensim-rcnl.jar\", \"Replace\");\nawait context.sync();\n", "ps1": "# Replace the placeholder path text \"{path to org-opensim-rcnl.jar}\" with\n# \"[PATH_TO_NMSM-CORE]\\gui\\org-opensim-rcnl.jar\" in the GUI installation\n# instructions paragraph (\"Run the command \"opensim64 --reload ...\").\n$d = $word.ActiveDocument\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"{path to org-opensim-rcnl.jar}\"\n$find.MatchCase = $true\n$find.MatchWildcards = $false\n$found = $find.Execute()\n\nif ($found) {\n    $range.Text = \"[PATH_TO_NMSM-CORE]\\gui\\org-opensim-rcnl.jar\"\n}\n"}
